$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ----- Row 3 -----
$ws.Range("A3").Value = 112213246
$ws.Range("B3").Value = 90332
$ws.Range("C3").Value = "Ovaliderad"
$ws.Range("D3").Value = "LC"
$ws.Range("E3").Value = 4769
$ws.Range("F3").Value = "Svavelriska"
$ws.Range("G3").Value = "Lactarius scrobiculatus"
$ws.Range("H3").Value = "(Scop.:Fr.) Fr."
$ws.Range("P3").Value = "Ol-olssvarttjärnen, Jmt"
$ws.Range("Q3").Value = 446878
$ws.Range("R3").Value = 7033464
$ws.Range("S3").Value = 10
$ws.Range("T3").Value = "Jämtland"
$ws.Range("U3").Value = "Krokom"
$ws.Range("V3").Value = "Jämtland"
$ws.Range("W3").Value = "Alsen"
$ws.Range("Y3").Value = "'2023-09-19"
$ws.Range("AA3").Value = "'2023-09-19"
$ws.Range("AD3").Value = $false
$ws.Range("AE3").Value = $false
$ws.Range("AG3").Value = $false
$ws.Range("AW3").Value = "Erik Lundmark"
$ws.Range("AX3").Value = "Erik Lundmark"

# ----- Row 4 -----
$ws.Range("A4").Value = 112213259
$ws.Range("B4").Value = 90687
$ws.Range("C4").Value = "Ovaliderad"
$ws.Range("D4").Value = "LC"
$ws.Range("E4").Value = 5964
$ws.Range("F4").Value = "Fjällig taggsvamp s.str."
$ws.Range("G4").Value = "Sarcodon imbricatus s.str."
$ws.Range("H4").Value = "(L.:Fr.) P.Karst."
$ws.Range("P4").Value = "Ol-olssvarttjärnen, Jmt"
$ws.Range("Q4").Value = 446861
$ws.Range("R4").Value = 7033454
$ws.Range("S4").Value = 10
$ws.Range("T4").Value = "Jämtland"
$ws.Range("U4").Value = "Krokom"
$ws.Range("V4").Value = "Jämtland"
$ws.Range("W4").Value = "Alsen"
$ws.Range("Y4").Value = "'2023-09-19"
$ws.Range("AA4").Value = "'2023-09-19"
$ws.Range("AD4").Value = $false
$ws.Range("AE4").Value = $false
$ws.Range("AG4").Value = $false
$ws.Range("AW4").Value = "Erik Lundmark"
$ws.Range("AX4").Value = "Erik Lundmark"
